# Insert a new record row at row 206 of the "Piña" price sheet.
# This shifts the existing rows 206:284 down to 207:285 and the new
# row 206 is populated with the newest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 206, pushing all
# subsequent rows (and their data/styles) down by one.
$ws.Rows("206:206").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A206").Value = 4
$ws.Range("B206").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C206").Value = 'Los Lagos'
$ws.Range("D206").Value = 44784
$ws.Range("E206").Value = 10
$ws.Range("F206").Value = 'Fruta'
$ws.Range("G206").Value = 100108
$ws.Range("H206").Value = 'Tropicales y subtropicales'
$ws.Range("I206").Value = 100108005
$ws.Range("J206").Value = 'Piña'
$ws.Range("K206").Value = 'Caramelo'
$ws.Range("L206").Value = 'Primera'
$ws.Range("M206").Value = 60
$ws.Range("N206").Value = 23000
$ws.Range("O206").Value = 23000
$ws.Range("P206").Value = 23000
$ws.Range("Q206").Value = '$/caja 12 unidades'
$ws.Range("R206").Value = 'Ecuador'
$ws.Range("S206").Value = 1917
$ws.Range("T206").Value = 12
